$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 25.15000000000049
$ws.Range("H2").Value = 0.5194690326159898
$ws.Range("I2").Value = 0.5194690326159898
$ws.Range("L2").Value = 3.33169892206371
$ws.Range("M2").Value = "[-4.751099215382242, 11.414497059509662]"
$ws.Range("N2").Value = 0.4108053593198946
$ws.Range("O2").Value = 0.4108053593198946
$ws.Range("P2").Value = -2.239053022378311
$ws.Range("Q2").Value = "[-5.3649219890132, 0.8868159442565782]"
$ws.Range("R2").Value = 0.1560305338093926
$ws.Range("S2").Value = 0.1560305338093926
$ws.Range("T2").Value = 11.04517346409555
$ws.Range("U2").Value = "[6.941910077527757, 15.148436850663336]"
$ws.Range("V2").Value = 0.000002234270654266624
$ws.Range("W2").Value = 0.000002234270654266624
$ws.Range("X2").Value = 8.962362362362537
$ws.Range("Y2").Value = -3.549699699699771
$ws.Range("Z2").Value = 21.47442442442485
$ws.Range("F3").Value = 25.15000000000049
$ws.Range("H3").Value = 0.2735214695801116
$ws.Range("I3").Value = 0.2735214695801116
$ws.Range("L3").Value = 4.731129689343151
$ws.Range("M3").Value = "[-2.7663521438115612, 12.228611522497863]"
$ws.Range("N3").Value = 0.2102745747590866
$ws.Range("O3").Value = 0.2102745747590866
$ws.Range("P3").Value = -1.924579283280234
$ws.Range("Q3").Value = "[-4.987553502095508, 1.1383949355350413]"
$ws.Range("R3").Value = 0.2121930280103734
$ws.Range("S3").Value = 0.2121930280103734
$ws.Range("T3").Value = 11.17140655081738
$ws.Range("U3").Value = "[6.999150012216258, 15.343663089418502]"
$ws.Range("V3").Value = 0.000002461931837505205
$ws.Range("W3").Value = 0.000002461931837505205
$ws.Range("X3").Value = 7.703603603603753
$ws.Range("Y3").Value = -4.556706706706802
$ws.Range("Z3").Value = 19.96391391391431
$ws.Range("F4").Value = 25.15000000000049
$ws.Range("H4").Value = 0.317806732551001
$ws.Range("I4").Value = 0.317806732551001
$ws.Range("L4").Value = 4.532541597450292
$ws.Range("M4").Value = "[-2.8861274541232564, 11.95121064902384]"
$ws.Range("N4").Value = 0.22488673179285
$ws.Range("O4").Value = 0.22488673179285
$ws.Range("P4").Value = -2.956053147521927
$ws.Range("Q4").Value = "[-6.075632639374855, 0.16352634433100022]"
$ws.Range("R4").Value = 0.06271437702969251
$ws.Range("S4").Value = 0.06271437702969251
$ws.Range("T4").Value = 10.59597601374022
$ws.Range("U4").Value = "[6.488361797187636, 14.703590230292804]"
$ws.Range("V4").Value = 0.000004783921115691925
$ws.Range("W4").Value = 0.000004783921115691925
$ws.Range("X4").Value = 11.83233233233256
$ws.Range("Y4").Value = -0.6545545545545668
$ws.Range("Z4").Value = 24.3192192192197
$ws.Range("F5").Value = 25.15000000000049
$ws.Range("H5").Value = 0.4496615052201114
$ws.Range("I5").Value = 0.4496615052201114
$ws.Range("L5").Value = 3.709071502080785
$ws.Range("M5").Value = "[-3.743302967994938, 11.161445972156509]"
$ws.Range("N5").Value = 0.3214979732991758
$ws.Range("O5").Value = 0.3214979732991758
$ws.Range("P5").Value = 2.987500521431735
$ws.Range("Q5").Value = "[-0.09434212172942313, 6.069343164592894]"
$ws.Range("R5").Value = 0.05712290205369719
$ws.Range("S5").Value = 0.05712290205369719
$ws.Range("T5").Value = 11.22971247206092
$ws.Range("U5").Value = "[7.163496094687764, 15.295928849434086]"
$ws.Range("V5").Value = 0.000001386843713202524
$ws.Range("W5").Value = 0.000001386843713202524
$ws.Range("X5").Value = 13.19179179179205
$ws.Range("Y5").Value = 0.8559559559559737
$ws.Range("Z5").Value = 25.52762762762813
$ws.Range("F6").Value = 25.15000000000049
$ws.Range("H6").Value = 0.3442461142271227
$ws.Range("I6").Value = 0.3442461142271227
$ws.Range("L6").Value = 4.310980175208291
$ws.Range("M6").Value = "[-3.42869751494376, 12.050657865360343]"
$ws.Range("N6").Value = 0.2678767434759077
$ws.Range("O6").Value = 0.2678767434759077
$ws.Range("P6").Value = -2.704474156243465
$ws.Range("Q6").Value = "[-5.798895748968547, 0.38994743648161645]"
$ws.Range("R6").Value = 0.0851537407081111
$ws.Range("S6").Value = 0.0851537407081111
$ws.Range("T6").Value = 10.42426327317792
$ws.Range("U6").Value = "[6.360824842374695, 14.487701703981154]"
$ws.Range("V6").Value = 0.000005266084393218406
$ws.Range("W6").Value = 0.000005266084393218406
$ws.Range("X6").Value = 10.82532532532553
$ws.Range("Y6").Value = -1.560860860860894
$ws.Range("Z6").Value = 23.21151151151197
$ws.Range("F7").Value = 25.15000000000049
$ws.Range("H7").Value = 0.4197457662377609
$ws.Range("I7").Value = 0.4197457662377609
$ws.Range("L7").Value = 4.226882212249629
$ws.Range("M7").Value = "[-4.095982394972248, 12.549746819471506]"
$ws.Range("N7").Value = 0.3118273587120888
$ws.Range("O7").Value = 0.3118273587120888
$ws.Range("P7").Value = 2.635289933641888
$ws.Range("Q7").Value = "[-0.49057903299300065, 5.761158900276778]"
$ws.Range("R7").Value = 0.09641241402637934
$ws.Range("S7").Value = 0.09641241402637934
$ws.Range("T7").Value = 12.01320100440797
$ws.Range("U7").Value = "[7.503137859235849, 16.523264149580093]"
$ws.Range("V7").Value = 0.000002706033358723658
$ws.Range("W7").Value = 0.000002706033358723658
$ws.Range("X7").Value = 14.60160160160189
$ws.Range("Y7").Value = 2.089539539539579
$ws.Range("Z7").Value = 27.11366366366419
$ws.Range("F8").Value = 25.15000000000049
$ws.Range("H8").Value = 0.2316990785006248
$ws.Range("I8").Value = 0.2316990785006248
$ws.Range("L8").Value = 4.77346801156221
$ws.Range("M8").Value = "[-2.320856288619762, 11.867792311744182]"
$ws.Range("N8").Value = 0.1821164331955174
$ws.Range("O8").Value = 0.1821164331955174
$ws.Range("P8").Value = 2.886868924920351
$ws.Range("Q8").Value = "[-0.19497371824080734, 5.968711568081509]"
$ws.Range("R8").Value = 0.06566536260877576
$ws.Range("S8").Value = 0.06566536260877576
$ws.Range("T8").Value = 10.60869251232044
$ws.Range("U8").Value = "[6.754676257495987, 14.462708767144889]"
$ws.Range("V8").Value = 0.000001475576598242156
$ws.Range("W8").Value = 0.000001475576598242156
$ws.Range("X8").Value = 13.59459459459486
$ws.Range("Y8").Value = 1.258758758758782
$ws.Range("Z8").Value = 25.93043043043094
$ws.Range("F9").Value = 23.77000000000028
$ws.Range("H9").Value = 0.1237440093088056
$ws.Range("I9").Value = 0.1237440093088056
$ws.Range("L9").Value = 6.333125706043914
$ws.Range("M9").Value = "[-1.3563539165579126, 14.02260532864574]"
$ws.Range("N9").Value = 0.1041043308499634
$ws.Range("O9").Value = 0.1041043308499634
$ws.Range("P9").Value = 1.314500229429964
$ws.Range("Q9").Value = "[-0.4842895582110387, 3.113290017070966]"
$ws.Range("R9").Value = 0.1480223613868901
$ws.Range("S9").Value = 0.1480223613868901
$ws.Range("T9").Value = 12.87854896025435
$ws.Range("U9").Value = "[8.585946401331457, 17.171151519177236]"
$ws.Range("V9").Value = 0.0000002698667567191393
$ws.Range("W9").Value = 0.0000002698667567191393
$ws.Range("X9").Value = 18.79709709709731
$ws.Range("Y9").Value = 11.99207207207221
$ws.Range("Z9").Value = 25.60212212212242
$ws.Range("F10").Value = 23.77000000000028
$ws.Range("H10").Value = 0.2849462290123931
$ws.Range("I10").Value = 0.2849462290123931
$ws.Range("L10").Value = 4.338477937518284
$ws.Range("M10").Value = "[-2.9935577125865294, 11.670513587623098]"
$ws.Range("N10").Value = 0.2395966435668619
$ws.Range("O10").Value = 0.2395966435668619
$ws.Range("P10").Value = 2.849132076228581
$ws.Range("Q10").Value = "[-0.2515789912784623, 5.949843143735625]"
$ws.Range("R10").Value = 0.07078551350211892
$ws.Range("S10").Value = 0.07078551350211892
$ws.Range("T10").Value = 8.956006847301119
$ws.Range("U10").Value = "[5.100313843024338, 12.811699851577899]"
$ws.Range("V10").Value = 0.00002662266564268378
$ws.Range("W10").Value = 0.00002662266564268378
$ws.Range("X10").Value = 12.99141141141156
$ws.Range("Y10").Value = 1.261071071071086
$ws.Range("Z10").Value = 24.72175175175204
$ws.Range("F11").Value = 23.77000000000028
$ws.Range("H11").Value = 0.07881951767834727
$ws.Range("I11").Value = 0.07881951767834727
$ws.Range("L11").Value = 6.133429961160727
$ws.Range("M11").Value = "[-0.9062926975631225, 13.173152619884577]"
$ws.Range("N11").Value = 0.08609907985371601
$ws.Range("O11").Value = 0.08609907985371601
$ws.Range("P11").Value = 2.006342455445734
$ws.Range("Q11").Value = "[0.19497371824080822, 3.817711192650659]"
$ws.Range("R11").Value = 0.03071641546777859
$ws.Range("S11").Value = 0.03071641546777859
$ws.Range("T11").Value = 10.6480182230082
$ws.Range("U11").Value = "[6.900402941392731, 14.395633504623671]"
$ws.Range("V11").Value = 0.0000008044008066310226
$ws.Range("W11").Value = 0.0000008044008066310226
$ws.Range("X11").Value = 16.17977977977997
$ws.Range("Y11").Value = 9.327167167167278
$ws.Range("Z11").Value = 23.03239239239266
$ws.Range("F12").Value = 23.77000000000028
$ws.Range("H12").Value = 0.05755932982514989
$ws.Range("I12").Value = 0.05755932982514989
$ws.Range("L12").Value = 6.709761425939946
$ws.Range("M12").Value = "[-0.13069300411091334, 13.550215855990805]"
$ws.Range("N12").Value = 0.05434922727349001
$ws.Range("O12").Value = 0.05434922727349001
$ws.Range("P12").Value = 1.754763464167272
$ws.Range("Q12").Value = "[0.2956053147521933, 3.213921613582351]"
$ws.Range("R12").Value = 0.01952032522537084
$ws.Range("S12").Value = 0.01952032522537084
$ws.Range("T12").Value = 12.71787502585404
$ws.Range("U12").Value = "[8.867043569889656, 16.56870648181842]"
$ws.Range("V12").Value = 0.00000003351567645637488
$ws.Range("W12").Value = 0.00000003351567645637488
$ws.Range("X12").Value = 17.13153153153173
$ws.Range("Y12").Value = 11.6113713713715
$ws.Range("Z12").Value = 22.65169169169195
$ws.Range("F13").Value = 23.77000000000028
$ws.Range("H13").Value = 0.1204065632322225
$ws.Range("I13").Value = 0.1204065632322225
$ws.Range("L13").Value = 5.492722245221623
$ws.Range("M13").Value = "[-1.7969828057589012, 12.782427296202147]"
$ws.Range("N13").Value = 0.13610799402365
$ws.Range("O13").Value = 0.13610799402365
$ws.Range("P13").Value = 2.446605690183042
$ws.Range("Q13").Value = "[-0.5786316799404627, 5.471843060306547]"
$ws.Range("R13").Value = 0.1103223997231608
$ws.Range("S13").Value = 0.1103223997231608
$ws.Range("T13").Value = 11.47679023792635
$ws.Range("U13").Value = "[7.765563218020004, 15.1880172578327]"
$ws.Range("V13").Value = 0.0000001428679019532808
$ws.Range("W13").Value = 0.0000001428679019532808
$ws.Range("X13").Value = 14.51421421421438
$ws.Range("Y13").Value = 3.069399399399432
$ws.Range("Z13").Value = 25.95902902902933
$ws.Range("F14").Value = 23.77000000000028
$ws.Range("H14").Value = 0.1433140936089161
$ws.Range("I14").Value = 0.1433140936089161
$ws.Range("L14").Value = 5.174968526775837
$ws.Range("M14").Value = "[-1.9780093322977645, 12.327946385849438]"
$ws.Range("N14").Value = 0.1520198100407577
$ws.Range("O14").Value = 0.1520198100407577
$ws.Range("P14").Value = 2.320816194543811
$ws.Range("Q14").Value = "[-0.798763297309117, 5.440395686396739]"
$ws.Range("R14").Value = 0.1410133093895243
$ws.Range("S14").Value = 0.1410133093895243
$ws.Range("T14").Value = 9.428542081932104
$ws.Range("U14").Value = "[5.772815096185443, 13.084269067678765]"
$ws.Range("V14").Value = 0.0000047994116159078
$ws.Range("W14").Value = 0.0000047994116159078
$ws.Range("X14").Value = 14.99009009009026
$ws.Range("Y14").Value = 3.188368368368405
$ws.Range("Z14").Value = 26.79181181181212
